$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Marking" row (row 11): marks per right answer and penalty per wrong answer
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Update the "Total" row (row 12): totals recomputed from the corrected marking scheme
$ws.Range("B12").Value = 84
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "78 / 112"
